# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Femacal de La Calera - Mango) right
# before the existing row 569, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 569:570 (existing data shifts down to 571:616).
$ws.Range("A569:A570").EntireRow.Insert()

# --- New row 569 ---
$ws.Range("A569").Value = 3
$ws.Range("B569").Value = "Femacal de La Calera"
$ws.Range("C569").Value = "Coquimbo"
$ws.Range("D569").Value = 45021
$ws.Range("E569").Value = 5
$ws.Range("F569").Value = "Fruta"
$ws.Range("G569").Value = 100108
$ws.Range("H569").Value = "Tropicales y subtropicales"
$ws.Range("I569").Value = 100108002
$ws.Range("J569").Value = "Mango"
$ws.Range("K569").Value = "Sin especificar"
$ws.Range("L569").Value = "Primera"
$ws.Range("M569").Value = 228
$ws.Range("N569").Value = 7000
$ws.Range("O569").Value = 7000
$ws.Range("P569").Value = 7000
$ws.Range("Q569").Value = "$/bandeja 4 kilos"
$ws.Range("R569").Value = "Perú"
$ws.Range("S569").Value = 1750
$ws.Range("T569").Value = 4

# --- New row 570 ---
$ws.Range("A570").Value = 3
$ws.Range("B570").Value = "Femacal de La Calera"
$ws.Range("C570").Value = "Coquimbo"
$ws.Range("D570").Value = 45021
$ws.Range("E570").Value = 5
$ws.Range("F570").Value = "Fruta"
$ws.Range("G570").Value = 100108
$ws.Range("H570").Value = "Tropicales y subtropicales"
$ws.Range("I570").Value = 100108002
$ws.Range("J570").Value = "Mango"
$ws.Range("K570").Value = "Sin especificar"
$ws.Range("L570").Value = "Segunda"
$ws.Range("M570").Value = 228
$ws.Range("N570").Value = 7000
$ws.Range("O570").Value = 7000
$ws.Range("P570").Value = 7000
$ws.Range("Q570").Value = "$/bandeja 4 kilos"
$ws.Range("R570").Value = "Perú"
$ws.Range("S570").Value = 1750
$ws.Range("T570").Value = 4
